# Cartoons workbook: add an English-title column (new column C), shifting the
# existing vote-matrix columns (old C:R) one place right to D:S.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column at C; everything that was C:R becomes D:S.
$ws.Columns("C").Insert()

# Match the new column's width to column B (poster filenames), same as the
# author did (cols min="2" max="3" width="27.33203125").
$ws.Columns("C").Width = $ws.Columns("B").Width

# English titles, row by row (this also matches the order in which the
# author must have typed them -- new shared strings 94..123 are exactly
# these titles in this row order, with the "cartoons-en" header being
# the very last new shared string, added after all the data rows).
$titles = @{
    4  = "Aladdin"
    5  = "Atlantis The Lost Empire"
    6  = "Brave"
    7  = "Cars"
    8  = "Coco"
    9  = "Despicable me"
    10 = "Encanto"
    11 = "Finding Nemo"
    12 = "Frozen"
    13 = "Happy Feet"
    14 = "Hotel Transylvania"
    15 = "How to Train Your Dragon"
    16 = "Ice Age"
    17 = "Incredibles"
    18 = "Kung Fu Panda"
    19 = "Lilo & Stitch"
    20 = "Madagascar"
    21 = "Megamind"
    22 = "Moana"
    23 = "Monsters Inc"
    24 = "Princess and the Frog"
    25 = "Puss in Boots"
    26 = "Rango"
    27 = "Ratatouille"
    28 = "Shrek"
    29 = "Spider Man - Into the Spider-Verse"
    30 = "Spirited Away"
    31 = "Tangled"
    32 = "The Boss Baby"
    33 = "The Croods"
    34 = "The Emperor's New Groove"
    35 = "The Lion King"
    36 = "Toy Story"
    37 = "Turbo"
    38 = "Up"
    39 = "WALL-E"
}

for ($r = 4; $r -le 39; $r++) {
    $ws.Cells.Item($r, 3).Value = $titles[$r]
}

# Header row (new column header, added last so it gets the final shared
# string slot) and the two metadata rows under it.
$ws.Range("C1").Value = "cartoons-en"
$ws.Range("C2").Value = "string"
$ws.Range("C3").Value = "meta"

# Author ended up with C1 selected.
$ws.Range("C1").Select() | Out-Null
